$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.904.39'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '3.390.52'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'571.13"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').Value = "'142.13"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.28%  '
$ws.Range('D7').Value = '3.390.49'
$ws.Range('E7').Value = '  -1.39%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.34%  '
$ws.Range('D10').Value = "'7.53"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').Value = "'0.395"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.32%  '
$ws.Range('D13').Value = '3.968.74'
$ws.Range('E14').Value = '  +1.94%  '
$ws.Range('D15').Value = "'28.26"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').Value = '3.389.60'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '60.984.90'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('E19').Value = '  -1.75%  '
$ws.Range('D20').Value = "'13.87"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').Value = "'8.97"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.55%  '
$ws.Range('D22').Value = "'384.64"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.83%  '
$ws.Range('D23').Value = "'0.558"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.44%  '
$ws.Range('D24').Value = "'74.41"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = "'0.0000118"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.76%  '
$ws.Range('D27').Value = '3.527.52'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('E28').Value = '  -1.19%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('D31').Value = "'7.98"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('E32').Value = '  -2.44%  '
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('D36').Value = "'6.99"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').Value = "'167.07"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').Value = '3.420.48'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('D39').Value = "'4.97"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.60%  '
$ws.Range('D40').Value = "'1.49"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.40%  '
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('D42').Value = "'27.32"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('D43').Value = "'0.781"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('D44').Value = "'1.00"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = "'42.15"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('D49').Value = '2.478.81'
$ws.Range('E49').Value = '  -4.65%  '
$ws.Range('D50').Value = "'6.83"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').Value = "'23.05"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.72%  '
